# The edit re-orders the observation records that live in worksheet rows
# 2-11 (each row is one species observation). The row that used to be
# row 6 is now row 2, the row that used to be row 11 is now row 3, etc.
# Columns C, I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AS, AT, AW, AX, AY are
# identical across all ten rows, so only the columns that actually carry
# per-row data (A, B, D, E, F, G, H, Q, R) need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")
$firstRow = 2
$lastRow = 11

# Snapshot the current ("before") values of the varying columns for every
# data row so the write-back below can't clobber a value it still needs to
# read (the permutation is cyclic, so naive in-place writes would corrupt
# data without this staging step).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $rowValues[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowValues
}

# old row number -> new row number, taken from the commit's row reshuffle.
$rowMap = @{
    2  = 8
    3  = 9
    4  = 10
    5  = 6
    6  = 2
    7  = 5
    8  = 11
    9  = 7
    10 = 4
    11 = 3
}

foreach ($oldRow in $rowMap.Keys) {
    $newRow = $rowMap[$oldRow]
    $rowValues = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $rowValues[$c]
    }
}
